$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2  = @{ B = 3.286832544864788;  C = 1.655778082260271; D = 0.7527432677738641; E = 0.4942365360607697; G = 6.189590430959694 }
    3  = @{ B = 0.1190320826869504; C = 1.655778082260271; D = 3.537761648806719;  E = 10.19245300693656;  G = 15.5050248206905 }
    4  = @{ B = 3.286832544864788;  C = 1.655778082260271; D = 3.537761648806719;  E = 0.4942365360607697; G = 8.974608811992548 }
    5  = @{ B = 0.2917716402565462; C = 1.655778082260271; D = 0.7527432677738641; E = 10.19245300693656;  G = 12.89274599722724 }
    6  = @{ B = 0.1190320826869504; C = 1.655778082260271; D = 0.7527432677738641; E = 1133.036916526867;  G = 1135.564469959588 }
    7  = @{ B = 0.2917716402565462; C = 0.306821227259698;  D = 0.7527432677738641; E = 10.19245300693656; G = 11.54378914222666 }
    8  = @{ B = 3.286832544864788;  C = 1.655778082260271; D = 3.537761648806719;  E = 10.19245300693656;  G = 18.67282528286833 }
    9  = @{ B = 3.286832544864788;  C = 1.655778082260271; D = 0.1494219747398047; E = 0.4942365360607697; G = 5.586269137925634 }
    10 = @{ B = 0.6606524410359556; C = 1.655778082260271; D = 6708.013860684405;  E = 1133.036916526867;  G = 7843.367207734568 }
    11 = @{ B = 3.286832544864788;  C = 117.745847958593;  D = 3.537761648806719;  E = 2195978.878461985;  G = 2196103.448904137 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
